$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "22.452.43"
Set-TextValue "E2" "  +0.18%  "

Set-TextValue "D3" "1.573.63"
Set-TextValue "E3" "  +0.68%  "

Set-TextValue "E4" "  +0.06%  "

Set-TextValue "E5" "  +0.05%  "

Set-TextValue "D6" "287.87"
Set-TextValue "E6" "  +0.56%  "

Set-TextValue "D7" "0.3721"
Set-TextValue "E7" "  +2.26%  "

Set-TextValue "D8" "47.54"
Set-TextValue "E8" "  -1.72%  "

Set-TextValue "D9" "0.3321"
Set-TextValue "E9" "  -0.62%  "

Set-TextValue "D10" "1.153"
Set-TextValue "E10" "  +2.20%  "

Set-TextValue "D11" "0.07525"
Set-TextValue "E11" "  +1.44%  "

Set-TextValue "E12" "  +0.10%  "

Set-TextValue "D13" "20.78"
Set-TextValue "E13" "  -0.24%  "

Set-TextValue "D14" "5.935"
Set-TextValue "E14" "  -0.04%  "

Set-TextValue "D15" "6.931"
Set-TextValue "E15" "  +0.67%  "

Set-TextValue "D16" "1.568.84"
Set-TextValue "E16" "  +0.38%  "

Set-TextValue "E17" "  +1.17%  "

Set-TextValue "D18" "88.32"
Set-TextValue "E18" "  -0.37%  "

Set-TextValue "D19" "0.06728"
Set-TextValue "E19" "  +0.42%  "

Set-TextValue "B20" "Dai"
Set-TextValue "C20" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D20" "1.000"
Set-TextValue "E20" "  +0.03%  "

Set-TextValue "B21" "Uniswap"
Set-TextValue "C21" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D21" "6.397"
Set-TextValue "E21" "  +0.93%  "

Set-TextValue "E22" "  +2.66%  "

Set-TextValue "D23" "12.01"
Set-TextValue "E23" "  +0.29%  "

Set-TextValue "D24" "22.435.14"
Set-TextValue "E24" "  +0.11%  "

Set-TextValue "D25" "2.388"
Set-TextValue "E25" "  -0.36%  "

Set-TextValue "D26" "2.622"
Set-TextValue "E26" "  +2.95%  "

Set-TextValue "D27" "150.82"
Set-TextValue "E27" "  +0.81%  "

Set-TextValue "D28" "19.64"
Set-TextValue "E28" "  +1.29%  "

Set-TextValue "D29" "4.939"
Set-TextValue "E29" "  -1.35%  "

Set-TextValue "D30" "125.29"

Set-TextValue "D31" "1.746.23"
Set-TextValue "E31" "  +0.50%  "

Set-TextValue "D32" "1.095"
Set-TextValue "E32" "  +3.19%  "

Set-TextValue "D33" "6.096"
Set-TextValue "E33" "  -0.55%  "

Set-TextValue "D34" "1.989"
Set-TextValue "E34" "  -0.43%  "

Set-TextValue "D35" "9.852"
Set-TextValue "E35" "  +2.58%  "

Set-TextValue "D36" "0.08340"
Set-TextValue "E36" "  +1.41%  "

Set-TextValue "E37" "  +2.55%  "

Set-TextValue "B38" "TrustWalletToken"
Set-TextValue "C38" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D38" "1.316"
Set-TextValue "E38" "  +1.06%  "

Set-TextValue "B39" "Algorand"
Set-TextValue "C39" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D39" "0.2234"
Set-TextValue "E39" "  +1.09%  "

Set-TextValue "D40" "0.06397"
Set-TextValue "E40" "  +0.02%  "

Set-TextValue "D41" "5.348"
Set-TextValue "E41" "  -0.12%  "

Set-TextValue "E42" "  +2.28%  "

Set-TextValue "E43" "  +3.16%  "

Set-TextValue "E44" "  +0.09%  "

Set-TextValue "D45" "14.04"
Set-TextValue "E45" "  +2.82%  "

Set-TextValue "D46" "0.6096"
Set-TextValue "E46" "  +5.97%  "

Set-TextValue "D47" "3.772"
Set-TextValue "E47" "  +0.25%  "

Set-TextValue "D48" "2.048"
Set-TextValue "E48" "  +1.74%  "

Set-TextValue "E49" "  +0.15%  "

Set-TextValue "E50" "  -0.23%  "

Set-TextValue "D51" "0.07203"
Set-TextValue "E51" "  -0.14%  "
